$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "purpose" column (G) for rows 2-63 was mislabeled "fullRNASeq";
# correct it to "spikein" for every data row.
$ws.Range("G2:G63").Value = "spikein"

# Reflect the selection left behind by the edit (G3:G63, active cell G3).
$ws.Range("G3:G63").Select()
